$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.489.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.96%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.849.50"
$ws.Range("D3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6291"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.18"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.30%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07531"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.2977"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.33%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "24.38"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.87%  "

$ws.Range("E12").Value = "  +0.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.906.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.40%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.008"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.84%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6869"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.59%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.74"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.28%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000009782"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.178.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.62%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.258"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.88%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "29.561.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.83%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "233.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.632"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.31%  "

$ws.Range("E25").Value = "  -0.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.77"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.42%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1394"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.441"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.37%  "

$ws.Range("E29").Value = "  -1.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.480"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05847"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.58%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.258"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.36%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.102"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.23%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.042"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.886"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.85%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.169"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.33%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.7194"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.57%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.585"
$ws.Range("D38").Style = "Normal"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.240.61"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.81%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.795"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.96%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01785"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.57%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9071"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.155"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.36%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.084.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.52%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "67.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.21%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.320"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.58%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.185"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.23%  "

$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4041"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.73%  "

$ws.Range("E51").Value = "  -1.18%  "
